$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.065483331680298
$ws.Range("B1").Value = 2.257797241210938
$ws.Range("C1").Value = 2.379379510879517
$ws.Range("D1").Value = 3.188571929931641
$ws.Range("E1").Value = 2.603131771087646
